$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K82").Value = 723
$ws.Range("I82").Value = 241
$ws.Range("M82").Value = -317
$ws.Range("H82").Value = 241
$ws.Range("K85").Value = 723
$ws.Range("I85").Value = 241
$ws.Range("M85").Value = 681
$ws.Range("H85").Value = 241
$ws.Range("K86").Value = 1862.1666
$ws.Range("I86").Value = 1862.1666
$ws.Range("M86").Value = -739.1666
$ws.Range("J86").Value = 2768
$ws.Range("N86").Value = -5014
$ws.Range("H86").Value = 2164.111
$ws.Range("L86").Value = 2768
$ws.Range("J88").Value = 50500
$ws.Range("H88").Value = 50500
$ws.Range("L88").Value = 50500
$ws.Range("N88").Value = -51312
$ws.Range("K89").Value = 9310.833000000001
$ws.Range("I89").Value = 1862.1666
$ws.Range("M89").Value = -3694.833000000001
$ws.Range("J89").Value = 2768
$ws.Range("N89").Value = -25072
$ws.Range("H89").Value = 2164.111
$ws.Range("L89").Value = 13840
$ws.Range("J91").Value = 50500
$ws.Range("H91").Value = 50500
$ws.Range("L91").Value = 50500
$ws.Range("N91").Value = -53308
$ws.Range("K107").Value = 678.2353000000001
$ws.Range("I107").Value = 678.2353000000001
$ws.Range("M107").Value = 1241.7647
$ws.Range("J107").Value = 1323.6666
$ws.Range("N107").Value = -5163.6666
$ws.Range("H107").Value = 846.6087
$ws.Range("L107").Value = 1323.6666
$ws.Range("K116").Value = 1821.7142
$ws.Range("I116").Value = 1821.7142
$ws.Range("M116").Value = 1620.2858
$ws.Range("J116").Value = 2440
$ws.Range("N116").Value = -9324
$ws.Range("H116").Value = 2079.3333
$ws.Range("L116").Value = 2440
$ws.Range("J133").Value = 0
$ws.Range("H133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("J136").Value = 49580
$ws.Range("H136").Value = 49580
$ws.Range("L136").Value = 49580
$ws.Range("N136").Value = -59780
$ws.Range("K137").Value = 18753970.5
$ws.Range("I137").Value = 6251323.5
$ws.Range("M137").Value = -18751420.5
$ws.Range("J137").Value = 2737.75
$ws.Range("N137").Value = -13313.25
$ws.Range("H137").Value = 2779887
$ws.Range("L137").Value = 8213.25
$ws.Range("J139").Value = 39797.273
$ws.Range("H139").Value = 39797.273
$ws.Range("L139").Value = 39797.273
$ws.Range("N139").Value = -50077.273
$ws.Range("J140").Value = 210780
$ws.Range("H140").Value = 210780
$ws.Range("L140").Value = 210780
$ws.Range("N140").Value = -221140

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K32").Value = 20094.5
$ws.Range("I32").Value = 20094.5
$ws.Range("M32").Value = -19807.5
$ws.Range("J32").Value = 29667.666
$ws.Range("N32").Value = -30241.666
$ws.Range("H32").Value = 21824.59
$ws.Range("L32").Value = 29667.666
$ws.Range("K45").Value = 1429.3334
$ws.Range("I45").Value = 1429.3334
$ws.Range("M45").Value = -1052.3334
$ws.Range("H45").Value = 1468
$ws.Range("J54").Value = 13000.5
$ws.Range("H54").Value = 13000.5
$ws.Range("L54").Value = 13000.5
$ws.Range("N54").Value = -14538.5
$ws.Range("K74").Value = 13932015
$ws.Range("I74").Value = 13932015
$ws.Range("M74").Value = -13931141
$ws.Range("H74").Value = 11179727
$ws.Range("K77").Value = 69660075
$ws.Range("I77").Value = 13932015
$ws.Range("M77").Value = -69655707
$ws.Range("H77").Value = 11179727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K31").Value = 4881.8
$ws.Range("I31").Value = 4881.8
$ws.Range("M31").Value = -4586.8
$ws.Range("J31").Value = 3536.7896
$ws.Range("N31").Value = -4126.7896
$ws.Range("H31").Value = 3817
$ws.Range("L31").Value = 3536.7896
$ws.Range("K34").Value = 4881.8
$ws.Range("I34").Value = 4881.8
$ws.Range("M34").Value = -4679.8
$ws.Range("J34").Value = 3536.7896
$ws.Range("N34").Value = -3940.7896
$ws.Range("H34").Value = 3817
$ws.Range("L34").Value = 3536.7896
$ws.Range("J133").Value = 35547.617
$ws.Range("H133").Value = 35547.617
$ws.Range("L133").Value = 35547.617
$ws.Range("N133").Value = -40607.617
$ws.Range("K134").Value = 3795.24
$ws.Range("I134").Value = 1265.08
$ws.Range("M134").Value = -1260.24
$ws.Range("J134").Value = 168516.5
$ws.Range("N134").Value = -510619.5
$ws.Range("H134").Value = 33636.324
$ws.Range("L134").Value = 505549.5
$ws.Range("J135").Value = 30000
$ws.Range("H135").Value = 30000
$ws.Range("L135").Value = 30000
$ws.Range("N135").Value = -40140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K23").Value = 393
$ws.Range("I23").Value = 131
$ws.Range("M23").Value = -158
$ws.Range("J23").Value = 184.42857
$ws.Range("N23").Value = -1023.28571
$ws.Range("H23").Value = 172.55556
$ws.Range("L23").Value = 553.28571
$ws.Range("K132").Value = 7282.8
$ws.Range("I132").Value = 809.2
$ws.Range("M132").Value = -4752.8
$ws.Range("J132").Value = 4079.8462
$ws.Range("N132").Value = -41778.6158
$ws.Range("H132").Value = 3171.3333
$ws.Range("L132").Value = 36718.6158
$ws.Range("K133").Value = 8257.5
$ws.Range("I133").Value = 2752.5
$ws.Range("M133").Value = -3197.5
$ws.Range("J133").Value = 7838.722
$ws.Range("N133").Value = -33636.166
$ws.Range("H133").Value = 5804.2334
$ws.Range("L133").Value = 23516.166
$ws.Range("K134").Value = 5421.1875
$ws.Range("I134").Value = 1807.0625
$ws.Range("M134").Value = -351.1875
$ws.Range("J134").Value = 7800
$ws.Range("N134").Value = -33540
$ws.Range("H134").Value = 3964.52
$ws.Range("L134").Value = 23400
$ws.Range("K136").Value = 8057.400000000001
$ws.Range("I136").Value = 2685.8
$ws.Range("M136").Value = -2957.400000000001
$ws.Range("J136").Value = 3647.5715
$ws.Range("N136").Value = -21142.7145
$ws.Range("H136").Value = 3246.8333
$ws.Range("L136").Value = 10942.7145
$ws.Range("K137").Value = 5250
$ws.Range("I137").Value = 1750
$ws.Range("M137").Value = -150
$ws.Range("J137").Value = 4000
$ws.Range("N137").Value = -22200
$ws.Range("H137").Value = 2200
$ws.Range("L137").Value = 12000
$ws.Range("K139").Value = 4754.85
$ws.Range("I139").Value = 1584.95
$ws.Range("M139").Value = 385.1499999999996
$ws.Range("J139").Value = 7583.7036
$ws.Range("N139").Value = -33031.11079999999
$ws.Range("H139").Value = 5031.0425
$ws.Range("L139").Value = 22751.1108
$ws.Range("K141").Value = 11000.0001
$ws.Range("I141").Value = 3666.6667
$ws.Range("M141").Value = -5820.000100000001
$ws.Range("J141").Value = 19640
$ws.Range("N141").Value = -69280
$ws.Range("H141").Value = 10927.272
$ws.Range("L141").Value = 58920

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("J80").Value = 4192.6
$ws.Range("N80").Value = -6188.6
$ws.Range("H80").Value = 4192.6
$ws.Range("L80").Value = 4192.6
$ws.Range("K83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("J83").Value = 4192.6
$ws.Range("N83").Value = -30947
$ws.Range("H83").Value = 4192.6
$ws.Range("L83").Value = 20963

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J2").Value = 2960.8
$ws.Range("H2").Value = 2960.8
$ws.Range("L2").Value = 2960.8
$ws.Range("N2").Value = -3184.8
$ws.Range("K22").Value = 717.2727
$ws.Range("I22").Value = 717.2727
$ws.Range("M22").Value = -422.2727
$ws.Range("J22").Value = 750.25
$ws.Range("N22").Value = -1340.25
$ws.Range("H22").Value = 731.1579
$ws.Range("L22").Value = 750.25
$ws.Range("K27").Value = 717.2727
$ws.Range("I27").Value = 717.2727
$ws.Range("M27").Value = -610.2727
$ws.Range("J27").Value = 750.25
$ws.Range("N27").Value = -964.25
$ws.Range("H27").Value = 731.1579
$ws.Range("L27").Value = 750.25
$ws.Range("K82").Value = 596.6667
$ws.Range("I82").Value = 596.6667
$ws.Range("M82").Value = -235.6667
$ws.Range("J82").Value = 1370
$ws.Range("N82").Value = -2092
$ws.Range("H82").Value = 1080
$ws.Range("L82").Value = 1370
$ws.Range("K85").Value = 596.6667
$ws.Range("I85").Value = 596.6667
$ws.Range("M85").Value = 651.3333
$ws.Range("J85").Value = 1370
$ws.Range("N85").Value = -3866
$ws.Range("H85").Value = 1080
$ws.Range("L85").Value = 1370

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K2").Value = 3000000
$ws.Range("I2").Value = 3000000
$ws.Range("M2").Value = -2999888
$ws.Range("J2").Value = 20750.75
$ws.Range("N2").Value = -20974.75
$ws.Range("H2").Value = 616600.6
$ws.Range("L2").Value = 20750.75
$ws.Range("J46").Value = 44979.6
$ws.Range("H46").Value = 44979.6
$ws.Range("L46").Value = 44979.6
$ws.Range("N46").Value = -45441.6
$ws.Range("J95").Value = 69931.336
$ws.Range("H95").Value = 69931.336
$ws.Range("L95").Value = 69931.336
$ws.Range("N95").Value = -75423.336
$ws.Range("K96").Value = 1233.3334
$ws.Range("I96").Value = 1233.3334
$ws.Range("M96").Value = 139.6666
$ws.Range("J96").Value = 1752
$ws.Range("N96").Value = -4498
$ws.Range("H96").Value = 1440.8
$ws.Range("L96").Value = 1752
$ws.Range("K107").Value = 1451.25
$ws.Range("I107").Value = 483.75
$ws.Range("M107").Value = 468.75
$ws.Range("J107").Value = 1144.1428
$ws.Range("N107").Value = -7272.428400000001
$ws.Range("H107").Value = 791.93335
$ws.Range("L107").Value = 3432.4284
$ws.Range("J134").Value = 44979.6
$ws.Range("H134").Value = 44979.6
$ws.Range("L134").Value = 134938.8
$ws.Range("N134").Value = -140008.8
